$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (weekly data refresh) ---
# Row 2
$ws.Cells.Item(2,4).Value = 44389
$ws.Cells.Item(2,10).Value = 120
$ws.Cells.Item(2,11).Value = 12000
$ws.Cells.Item(2,12).Value = 13000
$ws.Cells.Item(2,13).Value = 12500
$ws.Cells.Item(2,16).Value = 962

# Row 3
$ws.Cells.Item(3,4).Value = 45096
$ws.Cells.Item(3,10).Value = 750
$ws.Cells.Item(3,11).Value = 14000
$ws.Cells.Item(3,12).Value = 15000
$ws.Cells.Item(3,13).Value = 14600
$ws.Cells.Item(3,16).Value = 1123

# Row 4
$ws.Cells.Item(4,4).Value = 44972
$ws.Cells.Item(4,9).Value = "Primera"
$ws.Cells.Item(4,10).Value = 350
$ws.Cells.Item(4,11).Value = 17000
$ws.Cells.Item(4,12).Value = 18000
$ws.Cells.Item(4,13).Value = 17429
$ws.Cells.Item(4,14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(4,16).Value = 1162
$ws.Cells.Item(4,17).Value = 15

# Row 5
$ws.Cells.Item(5,4).Value = 44943
$ws.Cells.Item(5,9).Value = "Segunda"
$ws.Cells.Item(5,10).Value = 350
$ws.Cells.Item(5,11).Value = 14000
$ws.Cells.Item(5,12).Value = 15000
$ws.Cells.Item(5,13).Value = 14429
$ws.Cells.Item(5,16).Value = 1110

# Row 6
$ws.Cells.Item(6,4).Value = 45155
$ws.Cells.Item(6,10).Value = 300
$ws.Cells.Item(6,11).Value = 16000
$ws.Cells.Item(6,12).Value = 17000
$ws.Cells.Item(6,13).Value = 16500
$ws.Cells.Item(6,16).Value = 1269

# Row 7
$ws.Cells.Item(7,4).Value = 44918
$ws.Cells.Item(7,9).Value = "Segunda"
$ws.Cells.Item(7,10).Value = 200
$ws.Cells.Item(7,13).Value = 12750
$ws.Cells.Item(7,16).Value = 981

# Row 8
$ws.Cells.Item(8,4).Value = 44616
$ws.Cells.Item(8,10).Value = 120
$ws.Cells.Item(8,11).Value = 19000
$ws.Cells.Item(8,12).Value = 20000
$ws.Cells.Item(8,13).Value = 19500
$ws.Cells.Item(8,16).Value = 1500

# Row 9
$ws.Cells.Item(9,4).Value = 45028
$ws.Cells.Item(9,11).Value = 14000
$ws.Cells.Item(9,12).Value = 15000
$ws.Cells.Item(9,13).Value = 14500
$ws.Cells.Item(9,16).Value = 1115

# Row 10
$ws.Cells.Item(10,4).Value = 44469
$ws.Cells.Item(10,10).Value = 140
$ws.Cells.Item(10,11).Value = 13000
$ws.Cells.Item(10,12).Value = 14000
$ws.Cells.Item(10,13).Value = 13500
$ws.Cells.Item(10,16).Value = 1038

# Row 11
$ws.Cells.Item(11,4).Value = 44379
$ws.Cells.Item(11,10).Value = 120
$ws.Cells.Item(11,11).Value = 12000
$ws.Cells.Item(11,12).Value = 13000
$ws.Cells.Item(11,13).Value = 12667
$ws.Cells.Item(11,16).Value = 974

# Row 12
$ws.Cells.Item(12,4).Value = 45154
$ws.Cells.Item(12,10).Value = 250
$ws.Cells.Item(12,11).Value = 17000
$ws.Cells.Item(12,12).Value = 18000
$ws.Cells.Item(12,13).Value = 17500
$ws.Cells.Item(12,16).Value = 1346

# Row 13
$ws.Cells.Item(13,4).Value = 44406
$ws.Cells.Item(13,9).Value = "Primera"
$ws.Cells.Item(13,10).Value = 160
$ws.Cells.Item(13,11).Value = 17000
$ws.Cells.Item(13,12).Value = 18000
$ws.Cells.Item(13,13).Value = 17500
$ws.Cells.Item(13,16).Value = 1346

# Row 14
$ws.Cells.Item(14,4).Value = 44832
$ws.Cells.Item(14,10).Value = 100
$ws.Cells.Item(14,11).Value = 13000
$ws.Cells.Item(14,12).Value = 14000
$ws.Cells.Item(14,13).Value = 13500
$ws.Cells.Item(14,16).Value = 1038

# Row 16
$ws.Cells.Item(16,4).Value = 44592
$ws.Cells.Item(16,10).Value = 120
$ws.Cells.Item(16,11).Value = 12000
$ws.Cells.Item(16,12).Value = 13000
$ws.Cells.Item(16,13).Value = 12500
$ws.Cells.Item(16,16).Value = 962

# Row 17
$ws.Cells.Item(17,4).Value = 44988
$ws.Cells.Item(17,10).Value = 750
$ws.Cells.Item(17,11).Value = 17000
$ws.Cells.Item(17,12).Value = 18000
$ws.Cells.Item(17,13).Value = 17400
$ws.Cells.Item(17,16).Value = 1338

# Row 18
$ws.Cells.Item(18,4).Value = 45049
$ws.Cells.Item(18,10).Value = 300
$ws.Cells.Item(18,11).Value = 13000
$ws.Cells.Item(18,12).Value = 14000
$ws.Cells.Item(18,13).Value = 13500
$ws.Cells.Item(18,16).Value = 1038

# Row 20
$ws.Cells.Item(20,4).Value = 44580
$ws.Cells.Item(20,10).Value = 160
$ws.Cells.Item(20,11).Value = 11000
$ws.Cells.Item(20,12).Value = 12000
$ws.Cells.Item(20,13).Value = 11500
$ws.Cells.Item(20,16).Value = 885

# Row 21
$ws.Cells.Item(21,4).Value = 44159
$ws.Cells.Item(21,10).Value = 100
$ws.Cells.Item(21,11).Value = 23000
$ws.Cells.Item(21,12).Value = 24000
$ws.Cells.Item(21,13).Value = 23500
$ws.Cells.Item(21,16).Value = 1808

# Row 22
$ws.Cells.Item(22,4).Value = 44914
$ws.Cells.Item(22,10).Value = 100
$ws.Cells.Item(22,11).Value = 14000
$ws.Cells.Item(22,12).Value = 15000
$ws.Cells.Item(22,13).Value = 14400
$ws.Cells.Item(22,16).Value = 1108

# Row 23
$ws.Cells.Item(23,4).Value = 44910
$ws.Cells.Item(23,10).Value = 50
$ws.Cells.Item(23,11).Value = 14000
$ws.Cells.Item(23,12).Value = 15000
$ws.Cells.Item(23,13).Value = 14500
$ws.Cells.Item(23,16).Value = 1115

# Row 24
$ws.Cells.Item(24,4).Value = 45100
$ws.Cells.Item(24,10).Value = 200
$ws.Cells.Item(24,11).Value = 15000
$ws.Cells.Item(24,12).Value = 16000
$ws.Cells.Item(24,13).Value = 15500
$ws.Cells.Item(24,14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(24,16).Value = 1192
$ws.Cells.Item(24,17).Value = 13

# Row 25
$ws.Cells.Item(25,4).Value = 44764
$ws.Cells.Item(25,10).Value = 200
$ws.Cells.Item(25,11).Value = 12000
$ws.Cells.Item(25,12).Value = 13000
$ws.Cells.Item(25,13).Value = 12500
$ws.Cells.Item(25,16).Value = 962

# Row 26
$ws.Cells.Item(26,4).Value = 45140
$ws.Cells.Item(26,10).Value = 200
$ws.Cells.Item(26,11).Value = 16000
$ws.Cells.Item(26,12).Value = 17000
$ws.Cells.Item(26,13).Value = 16500
$ws.Cells.Item(26,16).Value = 1269

# Row 27
$ws.Cells.Item(27,4).Value = 45092
$ws.Cells.Item(27,10).Value = 600
$ws.Cells.Item(27,11).Value = 13000
$ws.Cells.Item(27,12).Value = 14000
$ws.Cells.Item(27,13).Value = 13500
$ws.Cells.Item(27,16).Value = 1038

# Row 28
$ws.Cells.Item(28,4).Value = 44320
$ws.Cells.Item(28,10).Value = 160
$ws.Cells.Item(28,11).Value = 19000
$ws.Cells.Item(28,12).Value = 20000
$ws.Cells.Item(28,13).Value = 19500
$ws.Cells.Item(28,16).Value = 1500

# Row 29
$ws.Cells.Item(29,4).Value = 45141
$ws.Cells.Item(29,10).Value = 400
$ws.Cells.Item(29,13).Value = 16550
$ws.Cells.Item(29,16).Value = 1273

# Row 31
$ws.Cells.Item(31,4).Value = 44893
$ws.Cells.Item(31,10).Value = 900
$ws.Cells.Item(31,11).Value = 13000
$ws.Cells.Item(31,12).Value = 14000
$ws.Cells.Item(31,13).Value = 13444
$ws.Cells.Item(31,16).Value = 1034

# --- Append new row 32 with the latest weekly price report ---
$ws.Cells.Item(32,1).Value = 1
$ws.Cells.Item(32,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(32,3).Value = "Arica y Parinacota"
$ws.Cells.Item(32,4).Value = 44890
$ws.Cells.Item(32,5).Value = 15
$ws.Cells.Item(32,6).Value = 100114007
$ws.Cells.Item(32,7).Value = "Jengibre"
$ws.Cells.Item(32,8).Value = "Sin especificar"
$ws.Cells.Item(32,9).Value = "Primera"
$ws.Cells.Item(32,10).Value = 300
$ws.Cells.Item(32,11).Value = 14000
$ws.Cells.Item(32,12).Value = 15000
$ws.Cells.Item(32,13).Value = 14500
$ws.Cells.Item(32,14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(32,15).Value = "Perú"
$ws.Cells.Item(32,16).Value = 1115
$ws.Cells.Item(32,17).Value = 13
$ws.Cells.Item(32,18).Value = "Hortaliza"

# Match the date number format used by the other rows in column D
$ws.Cells.Item(32,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

